$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 5942
$ws.Range("L3").Value = 6472
$ws.Range("L4").Value = 1592
$ws.Range("L6").Value = 5325
$ws.Range("L7").Value = 19714
$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L6").Value = 153
$ws.Range("L7").Value = 638
$ws.Range("L8").Value = 1300
$ws.Range("L14").Value = 99
$ws.Range("L20").Value = 498
$ws.Range("L29").Value = 1110
$ws.Range("L34").Value = 110
$ws.Range("L36").Value = 248
$ws.Range("L37").Value = 748
$ws.Range("L42").Value = 634
$ws.Range("L55").Value = 204
$ws.Range("L63").Value = 58
$ws.Range("L65").Value = 385
$ws.Range("L70").Value = 35
$ws.Range("L73").Value = 160
$ws.Range("L75").Value = 71
$ws.Range("L78").Value = 255
$ws.Range("L79").Value = 548
$ws.Range("L83").Value = 430
$ws.Range("L84").Value = 189
$ws.Range("L85").Value = 977
$ws.Range("L87").Value = 55
$ws.Range("L89").Value = 275
$ws.Range("L94").Value = 246
$ws.Range("L95").Value = 279
$ws.Range("L97").Value = 161
$ws.Range("L99").Value = 344
$ws.Range("L101").Value = 19714
$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("L2").Value = 41
$ws.Range("L7").Value = 99
$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L3").Value = 205
$ws.Range("L6").Value = 154
$ws.Range("L7").Value = 638
$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("L6").Value = 77
$ws.Range("L7").Value = 275
$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L3").Value = 405
$ws.Range("L7").Value = 977
$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L2").Value = 390
$ws.Range("L3").Value = 459
$ws.Range("L6").Value = 319
$ws.Range("L7").Value = 1300
$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("L6").Value = 94
$ws.Range("L7").Value = 430
$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("L3").Value = 88
$ws.Range("L6").Value = 65
$ws.Range("L7").Value = 279
$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L2").Value = 225
$ws.Range("L3").Value = 264
$ws.Range("L7").Value = 748
$ws = $wb.Worksheets.Item("New City")
$ws.Range("L3").Value = 126
$ws.Range("L7").Value = 385
$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L2").Value = 99
$ws.Range("L6").Value = 74
$ws.Range("L7").Value = 344
$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("L6").Value = 57
$ws.Range("L7").Value = 189
$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L3").Value = 432
$ws.Range("L4").Value = 61
$ws.Range("L6").Value = 272
$ws.Range("L7").Value = 1110
$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("L4").Value = 15
$ws.Range("L7").Value = 153
$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L2").Value = 173
$ws.Range("L7").Value = 634
$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("L2").Value = 67
$ws.Range("L3").Value = 86
$ws.Range("L4").Value = 28
$ws.Range("L7").Value = 255
$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("L6").Value = 55
$ws.Range("L7").Value = 204
$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L3").Value = 177
$ws.Range("L7").Value = 548
$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L2").Value = 156
$ws.Range("L3").Value = 173
$ws.Range("L7").Value = 498
$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("L4").Value = 20
$ws.Range("L7").Value = 248
$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("L2").Value = 38
$ws.Range("L7").Value = 110
$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("L2").Value = 59
$ws.Range("L3").Value = 59
$ws.Range("L7").Value = 246
$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("L2").Value = 58
$ws.Range("L7").Value = 160
$ws = $wb.Worksheets.Item("West Town")
$ws.Range("L3").Value = 34
$ws.Range("L7").Value = 161
$ws = $wb.Worksheets.Item("O'Hare")
$ws.Range("L4").Value = 2
$ws.Range("L7").Value = 35
$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("L2").Value = 31
$ws.Range("L6").Value = 9
$ws.Range("L7").Value = 71
$ws = $wb.Worksheets.Item("Ukrainian Village")
$ws.Range("L6").Value = 20
$ws.Range("L7").Value = 55
